$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54/55: scores + new comments in column E ---
$ws.Range("D54").Value = 4
$ws.Range("E54").Value = "Correct use of Axios"

$ws.Range("D55").Value = 4
$ws.Range("E55").Value = "Correct use of socket.io to complete two-way communication"

# --- Row 59: score ---
$ws.Range("D59").Value = 4

# --- Row 62/63: scores ---
$ws.Range("D62").Value = 4
$ws.Range("D63").Value = 4

# --- Row 69: "Have you followed the required work allocation?" -> answer Yes ---
$ws.Range("B69").Value = "Yes"

# --- Row 72: "What have they done?" per-member descriptions ---
$ws.Range("B72").Value = "socket.io, nodeJS server, the chat/annotation interface"
$ws.Range("C72").Value = "service worker, MongoDB"
$ws.Range("D72").Value = "IndexedDb, Axios communication, Swagger documentation"

# --- Row 73: "Percentage" per-member work split (33.3% each) ---
$ws.Range("B73").Value = 0.333
$ws.Range("B73").NumberFormat = "0.00%"

$ws.Range("C73").ClearFormats()
$ws.Range("C73").Value = 0.333
$ws.Range("C73").NumberFormat = "0.00%"
$ws.Range("C73").Font.Name = "ArialMT"
$ws.Range("C73").Font.Size = 11
$ws.Range("C73").Font.Color = 2171169

$ws.Range("D31").Copy()
$ws.Range("D73").PasteSpecial(-4122)
$ws.Range("D73").Value = 0.333

# --- Data validation: D55 no longer needs the 0-4 list validation ---
$ws.Range("D55").Validation.Delete()

# --- Update the active selection to match the latest edit location ---
$ws.Range("E54").Select()

$wb.Save()
